$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.575.01'
$ws.Range('E2').Value = '  +2.02%  '

$ws.Range('D3').Value = '1.885.97'
$ws.Range('E3').Value = '  +0.68%  '

$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.41%  '

$ws.Range('D5').Value = "'244.65"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.25%  '

$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.30%  '

$ws.Range('D7').Value = "'0.4917"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.61%  '

$ws.Range('D8').Value = "'0.2946"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.91%  '

$ws.Range('D9').Value = "'0.06758"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.48%  '

$ws.Range('D10').Value = '1.890.02'
$ws.Range('E10').Value = '  +0.90%  '

$ws.Range('D11').Value = "'17.24"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.92%  '

$ws.Range('D12').Value = "'0.07318"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.33%  '

$ws.Range('D13').Value = "'90.13"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.97%  '

$ws.Range('D14').Value = "'5.125"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +5.42%  '

$ws.Range('D15').Value = "'0.6732"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.84%  '

$ws.Range('D16').Value = '30.565.52'
$ws.Range('E16').Value = '  +2.04%  '

$ws.Range('D17').Value = "'0.000007890"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.31%  '

$ws.Range('D18').Value = "'1.002"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.40%  '

$ws.Range('D19').Value = "'13.28"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.68%  '

$ws.Range('D20').Value = '2.135.65'
$ws.Range('E20').Value = '  +1.03%  '

$ws.Range('D21').Value = "'1.002"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.43%  '

$ws.Range('D22').Value = "'4.853"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.44%  '

$ws.Range('D23').Value = "'181.73"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +33.25%  '

$ws.Range('D24').Value = "'6.024"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +7.83%  '

$ws.Range('D25').Value = "'9.314"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.93%  '

$ws.Range('D26').Value = "'155.19"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.84%  '

$ws.Range('D27').Value = "'18.57"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +9.55%  '

$ws.Range('D28').Value = "'1.908"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.69%  '

$ws.Range('D29').Value = "'1.392"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.79%  '

$ws.Range('D30').Value = "'4.335"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.73%  '

$ws.Range('D31').Value = "'0.08926"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.56%  '

$ws.Range('D32').Value = "'4.023"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.81%  '

$ws.Range('D33').Value = "'0.05204"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.35%  '

$ws.Range('D34').Value = "'0.7331"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.32%  '

$ws.Range('D35').Value = "'1.122"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.54%  '

$ws.Range('D36').Value = "'2.684"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.84%  '

$ws.Range('D37').Value = "'0.01863"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +9.47%  '

$ws.Range('D38').Value = "'2.659"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.82%  '

$ws.Range('D39').Value = "'2.144"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.45%  '

$ws.Range('D40').Value = "'0.9361"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.87%  '

$ws.Range('D41').Value = "'0.4489"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +6.89%  '

$ws.Range('D42').Value = "'105.71"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.68%  '

$ws.Range('D43').Value = "'5.784"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.65%  '

$ws.Range('D44').Value = "'1.002"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.53%  '

$ws.Range('D45').Value = "'7.621"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.13%  '

$ws.Range('D46').Value = "'0.1340"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +6.96%  '

$ws.Range('D47').Value = "'0.05859"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.88%  '

$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = "'0.4020"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +8.31%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'8.558"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.90%  '

$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = "'33.34"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.22%  '

$ws.Range('D51').Value = "'1.397"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.96%  '
